$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.869.07'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = '2.032.59'
$ws.Range("E3").Value = '  -0.89%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.607'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.17'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.82%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.381'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.02%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0812'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("E11").Value = '  +1.06%  '

$ws.Range("D12").Value = '2.335.69'
$ws.Range("E12").Value = '  -0.74%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.55'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.67%  '

$ws.Range("E14").Value = '  +2.62%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.760'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.17%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.16'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.36%  '

$ws.Range("D17").Value = '2.038.74'
$ws.Range("E17").Value = '  +0.41%  '

$ws.Range("D18").Value = '37.836.57'
$ws.Range("E18").Value = '  -0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.87'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.26%  '

$ws.Range("D21").Value = '0.0₃0826'
$ws.Range("E21").Value = '  -1.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("E24").Value = '  -0.36%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.89%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.27%  '

$ws.Range("E28").Value = '  -3.87%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.57%  '

$ws.Range("E30").Value = '  -4.14%  '

$ws.Range("E31").Value = '  +0.96%  '

$ws.Range("E32").Value = '  +4.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.67%  '

$ws.Range("E34").Value = '  -0.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.56%  '

$ws.Range("E37").Value = '  -2.60%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.25'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.23%  '

$ws.Range("E39").Value = '  +0.14%  '

$ws.Range("D40").Value = '1.524.21'
$ws.Range("E40").Value = '  +2.49%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.22'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.04%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0217'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.87%  '

$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("E45").Value = '  -0.85%  '

$ws.Range("E46").Value = '  -1.78%  '

$ws.Range("E47").Value = '  -3.07%  '

$ws.Range("E48").Value = '  -0.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.96'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.14%  '

$ws.Range("E50").Value = '  +0.59%  '

$ws.Range("D51").Value = '2.224.62'
$ws.Range("E51").Value = '  -0.78%  '
